$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "UI" column (G) as "Done" for rows 9, 10, 11
$ws.Range("G9").Value = "Done"
$ws.Range("G10").Value = "Done"
$ws.Range("G11").Value = "Done"

# Update the selection to D13:D17 with active cell D13
$ws.Range("D13:D17").Select()
